$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Insert new columns C-F (CheckInDate, CheckOutDate, StayDuration, ChargesOfStay)
# ahead of the existing "OrdersOfUser" column (now G), and append "TotalBill" (H).
$ws.Range("C1").Value = "CheckInDate"
$ws.Range("D1").Value = "CheckOutDate"
$ws.Range("E1").Value = "StayDuration"
$ws.Range("F1").Value = "ChargesOfStay"
$ws.Range("G1").Value = "OrdersOfUser"
$ws.Range("H1").Value = "TotalBill"

# Force the date-like columns to be stored as plain text (not auto-converted
# to Excel date serials) and allow decimals to stay exactly as typed.
$ws.Range("C2:D4").NumberFormat = "@"
$ws.Range("E2:E4").NumberFormat = "@"

# --- Row 2: Abhishek Shigavan ---
$ws.Range("B2").Value = "Abhi9876"
$ws.Range("C2").Value = "10/04/2021"
$ws.Range("D2").Value = "10/04/2021"
$ws.Range("E2").Value = "0 Day 0 Hours"
$ws.Range("F2").Value = 590
$ws.Range("G2").Value = "1. Name Of Item : Breakfast, Quantity : 2, Price : 150, TotalAmount : 354 | 2. Name Of Item : Bedsheet, Quantity : 2, Price : 300, TotalAmount : 708 | 3. Name Of Item : Dinner, Quantity : 2, Price : 300, TotalAmount : 708 | "
$ws.Range("H2").Value = 2360

# --- Row 3: Peter Parkar ---
$ws.Range("B3").Value = "Pete9988"
$ws.Range("C3").Value = "10/04/2021"
$ws.Range("D3").Value = "10/04/2021"
$ws.Range("E3").Value = "0 Day 0 Hours"
$ws.Range("F3").Value = 590
$ws.Range("G3").Value = "1. Name Of Item : Tea, Quantity : 4, Price : 25, TotalAmount : 118 | 2. Name Of Item : Soap, Quantity : 1, Price : 40, TotalAmount : 47.2 | "
$ws.Range("H3").Value = 755.2

# --- Row 4: Rajat Sawarkar (new row) ---
$ws.Range("A4").Value = "Rajat Sawarkar"
$ws.Range("B4").Value = "Raja9899"
$ws.Range("C4").Value = "10/04/2021"
$ws.Range("D4").Value = "10/04/2021"
$ws.Range("E4").Value = "0 Day 0 Hours"
$ws.Range("F4").Value = 590
$ws.Range("G4").Value = "1. Name Of Item : Breakfast, Quantity : 3, Price : 150, TotalAmount : 531 | 2. Name Of Item : Blanket, Quantity : 2, Price : 400, TotalAmount : 944 | "
$ws.Range("H4").Value = 2065
